$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "version / round" counter in A1 (24 -> 25)
$ws.Range("A1").Value = 25

# Append a new vocabulary row (row 26) with the new lesson's words
$ws.Range("A26").Value = "lek 20"
$ws.Range("B26").Value = "test"
$ws.Range("C26").Value = "zum Beispiel"
$ws.Range("D26").Value = "الالللليب"
$ws.Range("E26").Value = "سيبسيبسيب"

# Match the formatting used by the rest of column A (centered style)
$ws.Range("A25").Copy()
$ws.Range("A26").PasteSpecial(-4122)
